$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so that plain numeric-looking
# strings (e.g. "65.17") are not auto-converted to numbers by Excel, then clear
# the format back off again so cell styling is left untouched.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.128.02"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "1.658.09"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "215.34"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "1.892.09"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("D13").Value = "1.664.83"
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "65.17"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "27.111.98"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "238.49"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").Value = "7.93"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("E23").Value = "  +5.01%  "
$ws.Range("D24").Value = "9.25"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "15.86"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "1.524.08"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +10.00%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "0.579"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").Value = "0.891"
$ws.Range("E38").Value = "  +8.37%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "66.48"
$ws.Range("E42").Value = "  +9.66%  "
$ws.Range("D43").Value = "2.26"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").Value = "1.799.55"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "0.779"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "90.11"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("D50").Value = "0.0506"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "0.0978"
$ws.Range("E51").Value = "  +3.08%  "

$dataRange.ClearFormats()

